$wb = $excel.ActiveWorkbook

# Sheet "P_valores"
$ws = $wb.Worksheets.Item("P_valores")
$ws.Range("C2").Value = 0.9681783791380707
$ws.Range("D2").Value = 0.558892892700422
$ws.Range("E2").Value = 0.9489550238307864
$ws.Range("F2").Value = 0.5007263934923354

$ws.Range("B3").Value = 0.9681783791380707
$ws.Range("D3").Value = 0.5569497369304761
$ws.Range("E3").Value = 0.8880657230589044
$ws.Range("F3").Value = 0.5441174782034848

$ws.Range("B4").Value = 0.558892892700422
$ws.Range("C4").Value = 0.5569497369304761
$ws.Range("E4").Value = 0.6028125286686925
$ws.Range("F4").Value = 0.1767662957478791

$ws.Range("B5").Value = 0.9489550238307864
$ws.Range("C5").Value = 0.8880657230589044
$ws.Range("D5").Value = 0.6028125286686925
$ws.Range("F5").Value = 0.5879119629688141

$ws.Range("B6").Value = 0.5007263934923354
$ws.Range("C6").Value = 0.5441174782034848
$ws.Range("D6").Value = 0.1767662957478791
$ws.Range("E6").Value = 0.5879119629688141

# Sheet "Estadisticos_DM"
$ws2 = $wb.Worksheets.Item("Estadisticos_DM")
$ws2.Range("C2").Value = 0.04061241660111244
$ws2.Range("D2").Value = 0.5987667449857927
$ws2.Range("E2").Value = -0.06517648862407824
$ws2.Range("F2").Value = -0.6912245288699611

$ws2.Range("B3").Value = -0.04061241660111244
$ws2.Range("D3").Value = 0.6017665082094732
$ws2.Range("E3").Value = -0.1433382472837913
$ws2.Range("F3").Value = -0.6217209777829579

$ws2.Range("B4").Value = -0.5987667449857927
$ws2.Range("C4").Value = -0.6017665082094732
$ws2.Range("E4").Value = -0.5323767987562255
$ws2.Range("F4").Value = -1.422544808300594

$ws2.Range("B5").Value = 0.06517648862407824
$ws2.Range("C5").Value = 0.1433382472837913
$ws2.Range("D5").Value = 0.5323767987562255
$ws2.Range("F5").Value = -0.5546112447853006

$ws2.Range("B6").Value = 0.6912245288699611
$ws2.Range("C6").Value = 0.6217209777829579
$ws2.Range("D6").Value = 1.422544808300594
$ws2.Range("E6").Value = 0.5546112447853006
